# Swap the full record (columns B:AC) between pairs of rows.
# Column A (the running index) is left untouched on each row.
# This matches the source diff, where the two fixtures listed
# consecutively in the sheet had their data rows exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(17, 18),
    @(27, 28),
    @(84, 85),
    @(88, 89),
    @(177, 178)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
